$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (scheduled GitHub Actions refresh)
# Leading apostrophe forces Excel to store values as text, preserving
# formatting such as trailing zeros and thousands separators written as dots.

$ws.Range("D2").Value = "'59.893.37"
$ws.Range("E2").Value = "'  +1.41%  "
$ws.Range("D3").Value = "'2.665.71"
$ws.Range("E3").Value = "'  +2.89%  "
$ws.Range("E4").Value = "'  -0.07%  "
$ws.Range("D5").Value = "'535.82"
$ws.Range("E5").Value = "'  +1.03%  "
$ws.Range("D6").Value = "'145.44"
$ws.Range("E6").Value = "'  +4.08%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "'  +0.19%  "
$ws.Range("E8").Value = "'  +1.52%  "
$ws.Range("D9").Value = "'2.663.75"
$ws.Range("E9").Value = "'  +2.40%  "
$ws.Range("D10").Value = "'6.62"
$ws.Range("E10").Value = "'  +2.54%  "
$ws.Range("E11").Value = "'  +2.22%  "
$ws.Range("D12").Value = "'0.338"
$ws.Range("E12").Value = "'  +1.52%  "
$ws.Range("E13").Value = "'  -1.56%  "
$ws.Range("D14").Value = "'3.119.42"
$ws.Range("E14").Value = "'  +2.31%  "
$ws.Range("D15").Value = "'59.827.10"
$ws.Range("E15").Value = "'  +1.38%  "
$ws.Range("D16").Value = "'21.11"
$ws.Range("E16").Value = "'  +3.13%  "
$ws.Range("D17").Value = "'2.683.35"
$ws.Range("E17").Value = "'  +2.44%  "
$ws.Range("D18").Value = "'0.0000135"
$ws.Range("E18").Value = "'  +1.44%  "
$ws.Range("D19").Value = "'344.07"
$ws.Range("E19").Value = "'  -0.95%  "
$ws.Range("D20").Value = "'4.44"
$ws.Range("E20").Value = "'  +2.38%  "
$ws.Range("D21").Value = "'10.28"
$ws.Range("E21").Value = "'  +1.81%  "
$ws.Range("D22").Value = "'6.38"
$ws.Range("E22").Value = "'  -0.83%  "
$ws.Range("E23").Value = "'  -0.08%  "
$ws.Range("D24").Value = "'67.82"
$ws.Range("E24").Value = "'  +0.34%  "
$ws.Range("D25").Value = "'0.415"
$ws.Range("E25").Value = "'  +2.38%  "
$ws.Range("D26").Value = "'0.166"
$ws.Range("E26").Value = "'  -0.46%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "'  +0.29%  "
$ws.Range("D28").Value = "'7.29"
$ws.Range("E28").Value = "'  +2.54%  "
$ws.Range("D29").Value = "'0.0₃0753"
$ws.Range("E29").Value = "'  +4.12%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "'  +0.02%  "
$ws.Range("E31").Value = "'  +2.65%  "
$ws.Range("D32").Value = "'5.94"
$ws.Range("E32").Value = "'  +1.02%  "
$ws.Range("D33").Value = "'19.08"
$ws.Range("E33").Value = "'  +1.79%  "
$ws.Range("D34").Value = "'150.44"
$ws.Range("E34").Value = "'  +0.85%  "
$ws.Range("D35").Value = "'4.03"
$ws.Range("E35").Value = "'  +1.59%  "
$ws.Range("E36").Value = "'  +2.45%  "
$ws.Range("D37").Value = "'1.47"
$ws.Range("E37").Value = "'  -0.21%  "
$ws.Range("D38").Value = "'0.840"
$ws.Range("E38").Value = "'  +1.43%  "
$ws.Range("D39").Value = "'0.824"
$ws.Range("E39").Value = "'  +0.53%  "
$ws.Range("D40").Value = "'289.62"
$ws.Range("E40").Value = "'  +7.64%  "
$ws.Range("D41").Value = "'3.60"
$ws.Range("E41").Value = "'  +1.90%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "'  +0.18%  "
$ws.Range("D43").Value = "'0.604"
$ws.Range("E43").Value = "'  +1.37%  "
$ws.Range("B44").Value = "'Hedera"
$ws.Range("C44").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").Value = "'0.0540"
$ws.Range("E44").Value = "'  +4.24%  "
$ws.Range("B45").Value = "'WhiteBITCoin"
$ws.Range("C45").Value = "'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").Value = "'10.73"
$ws.Range("E45").Value = "'  -0.09%  "
$ws.Range("B46").Value = "'Stellar"
$ws.Range("C46").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.0958"
$ws.Range("E46").Value = "'  +0.02%  "
$ws.Range("D47").Value = "'1.970.00"
$ws.Range("E47").Value = "'  +0.53%  "
$ws.Range("D48").Value = "'0.0225"
$ws.Range("E48").Value = "'  +1.83%  "
$ws.Range("D49").Value = "'4.58"
$ws.Range("E49").Value = "'  -1.03%  "
$ws.Range("D50").Value = "'18.46"
$ws.Range("E50").Value = "'  +1.02%  "
$ws.Range("D51").Value = "'110.01"
$ws.Range("E51").Value = "'  -1.00%  "
